# Apply the "gh-pages output regenerated" update to 广州-漫展信息.xlsx
# (想去人数 / F column headcounts bump up, one 展览 listing flips from
#  "已售罄" to "暂时售罄" and therefore re-enters the "全部类型" aggregate
#  sheet in place of the expired "广州·火影only" row.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 1504
$ws1.Range("F6").Value  = 736
$ws1.Range("F7").Value  = 44
$ws1.Range("F8").Value  = 669
$ws1.Range("F11").Value = 1393
$ws1.Range("F12").Value = 36600
$ws1.Range("G12").Value = "暂时售罄"
$ws1.Range("F13").Value = 7270
$ws1.Range("F15").Value = 383
$ws1.Range("F19").Value = 117
$ws1.Range("F20").Value = 430
$ws1.Range("F23").Value = 124
$ws1.Range("F24").Value = 829
$ws1.Range("F30").Value = 232
$ws1.Range("F35").Value = 770

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F6").Value  = 297
$ws2.Range("F8").Value  = 4
$ws2.Range("F10").Value = 5

# ---------------------------------------------------------------------
# Sheet "本地生活" (local-life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F3").Value = 369

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types combined)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 369
$ws4.Range("F5").Value  = 1504
$ws4.Range("F9").Value  = 736
$ws4.Range("F10").Value = 44
$ws4.Range("F11").Value = 669

# Row 13 used to be the (now stale) "广州·火影only" listing; it is
# replaced wholesale with the "萤火虫动漫游戏嘉年华" listing, which has
# just flipped out of "已售罄" and so re-qualifies for this sheet.
# Column B holds a literal text date ("2024-07-19"), not a real Excel
# date value, so force a text number format before assigning it to
# stop the COM layer auto-converting it to a date serial, then put the
# format back.
$cellB13 = $ws4.Cells.Item(13, 2)
$cellB13.NumberFormat = "@"
$cellB13.Value = "2024-07-19"
$cellB13.NumberFormat = "General"

$ws4.Cells.Item(13, 3).Value = "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园"
$ws4.Cells.Item(13, 4).Value = "新港东路1000号 保利世贸博览馆"
$ws4.Cells.Item(13, 5).Value = "2024.07.19 09:00-07.22 17:00"
$ws4.Cells.Item(13, 6).Value = 36600
$ws4.Cells.Item(13, 7).Value = 0
$ws4.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87210"
$ws4.Cells.Item(13, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg"

$ws4.Range("F15").Value = 297
$ws4.Range("F17").Value = 5
$ws4.Range("F19").Value = 7270
$ws4.Range("F20").Value = 383
$ws4.Range("F25").Value = 117
$ws4.Range("F26").Value = 430
$ws4.Range("F31").Value = 124
$ws4.Range("F32").Value = 829
$ws4.Range("F38").Value = 232
